$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking values
# (e.g. "333.07") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Updated cryptos list values (Price in column D, Volume(1h) in column E)
$ws.Range("D2").Value = '27.551.81'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '1.844.73'
$ws.Range("E3").Value = '  -2.40%  '
$ws.Range("E4").Value = '  -1.07%  '
$ws.Range("D5").Value = '333.07'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("D7").Value = '0.4640'
$ws.Range("E7").Value = '  -1.39%  '
$ws.Range("D8").Value = '0.3855'
$ws.Range("E8").Value = '  -1.97%  '
$ws.Range("D9").Value = '46.22'
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("D10").Value = '0.07906'
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("D11").Value = '0.9934'
$ws.Range("E11").Value = '  -2.23%  '
$ws.Range("D12").Value = '21.46'
$ws.Range("E12").Value = '  -1.63%  '
$ws.Range("D13").Value = '1.842.76'
$ws.Range("E13").Value = '  -2.53%  '
$ws.Range("D14").Value = '5.913'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").Value = '7.106'
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("D17").Value = '88.86'
$ws.Range("E17").Value = '  +1.39%  '
$ws.Range("D18").Value = '0.06647'
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("D19").Value = '0.00001034'
$ws.Range("E19").Value = '  -1.53%  '
$ws.Range("D20").Value = '17.05'
$ws.Range("E20").Value = '  -0.63%  '
$ws.Range("E21").Value = '  -1.04%  '
$ws.Range("D22").Value = '27.563.47'
$ws.Range("E22").Value = '  -1.47%  '
$ws.Range("D23").Value = '5.380'
$ws.Range("E23").Value = '  -2.15%  '
$ws.Range("D24").Value = '10.90'
$ws.Range("E24").Value = '  -0.74%  '
$ws.Range("D25").Value = '2.302'
$ws.Range("E25").Value = '  -2.76%  '
$ws.Range("D26").Value = '158.16'
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("D27").Value = '19.51'
$ws.Range("E27").Value = '  -2.56%  '
$ws.Range("D28").Value = '2.098'
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("D29").Value = '5.399'
$ws.Range("E29").Value = '  -2.04%  '
$ws.Range("D30").Value = '119.82'
$ws.Range("E30").Value = '  -1.29%  '
$ws.Range("D31").Value = '0.9755'
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("D32").Value = '0.09407'
$ws.Range("E32").Value = '  -1.75%  '
$ws.Range("D33").Value = '3.586'
$ws.Range("E33").Value = '  -1.75%  '
$ws.Range("D34").Value = '5.285'
$ws.Range("E34").Value = '  -1.27%  '
$ws.Range("D35").Value = '1.339'
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("D36").Value = '0.06024'
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("D37").Value = '0.02227'
$ws.Range("E37").Value = '  -1.09%  '
$ws.Range("D38").Value = '8.309'
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("D39").Value = '1.182'
$ws.Range("E39").Value = '  -3.00%  '
$ws.Range("D40").Value = '0.5886'
$ws.Range("E40").Value = '  -1.30%  '
$ws.Range("D41").Value = '0.1862'
$ws.Range("E41").Value = '  -2.00%  '
$ws.Range("D42").Value = '10.29'
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").Value = '1.243'
$ws.Range("E43").Value = '  -2.05%  '
$ws.Range("D44").Value = '0.5574'
$ws.Range("E44").Value = '  -2.00%  '
$ws.Range("D45").Value = '12.12'
$ws.Range("E45").Value = '  -1.33%  '
$ws.Range("D46").Value = '1.901'
$ws.Range("E46").Value = '  -2.31%  '
$ws.Range("D47").Value = '0.06686'
$ws.Range("E47").Value = '  -2.54%  '
$ws.Range("D48").Value = '110.76'
$ws.Range("E48").Value = '  -2.87%  '
$ws.Range("D49").Value = '1.052'
$ws.Range("E49").Value = '  -1.60%  '
$ws.Range("D50").Value = '1.003'
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("D51").Value = '70.02'
$ws.Range("E51").Value = '  -1.51%  '

# Restore default styling on column D so no residual number-format style remains
$ws.Range("D2:D51").Style = "Normal"
